$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 30 (pushes existing rows 30-34 down to 31-35)
$ws.Rows.Item(30).Insert()

# Fill in the new row's content
$ws.Range("A30").Value2 = "JPA"
$ws.Range("B30").Value2 = "Quick facts"

$factsText = @"
# JDBC rely on SQL, while new approach like JPA, Hibernate, EclipseLink, are more user-friendly.
# JPA base on JDBC in bg.
# JPA is the API providers. Such as EclipseLink and Hibernate define the concrete impl for the JPA specification.
# 2 ways to specify the ORM impl, the orm.xml or the annotation (@Entity, @Table, @Id, @Transient, @OnetoMany, @ManytoMany)
# Dependency for a simple JPA project can be: 
  - mysql-connector 
  - hibernate jpa
  - hibernate core
  - hibernate entity manager
"@
$factsText = $factsText.TrimEnd("`n")

$ws.Range("C30").Value2 = $factsText

# Bold the "orm.xml" substring within the cell's rich text
$boldStart = 276
$boldLength = 7
$chars = $ws.Range("C30").Characters($boldStart, $boldLength)
$chars.Font.Bold = $true

# Match the new row's height with its sibling rows
$ws.Rows.Item(30).RowHeight = 33

# Restore the originally selected cell
$ws.Range("C28").Select()
